$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.746.30"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.484.64"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'587.22"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'175.35"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'0.145"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'4.96"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "2.937.63"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'25.29"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "67.629.82"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "2.480.87"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'10.81"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'7.41"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "'346.80"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'70.73"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'4.19"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'1.69"
$ws.Range("E25").Value = "  -6.66%  "
$ws.Range("D26").Value = "'8.86"
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "0.0₃0893"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("D30").Value = "'496.85"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").Value = "'7.73"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'164.62"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("E36").Value = "  +1.85%  "
$ws.Range("D37").Value = "'18.63"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'147.92"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "'0.576"
$ws.Range("E51").Value = "  -1.38%  "
